$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 400
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("H62").Value = 8686.429
$ws.Range("I62").Value = 3701.25
$ws.Range("J62").Value = 15333.333
$ws.Range("K62").Value = 3701.25
$ws.Range("L62").Value = 15333.333
$ws.Range("M62").Value = -3077.25
$ws.Range("N62").Value = -16581.333
$ws.Range("H65").Value = 8686.429
$ws.Range("I65").Value = 3701.25
$ws.Range("J65").Value = 15333.333
$ws.Range("K65").Value = 18506.25
$ws.Range("L65").Value = 76666.66500000001
$ws.Range("M65").Value = -15386.25
$ws.Range("N65").Value = -82906.66500000001
$ws.Range("H132").Value = 2351.9167
$ws.Range("I132").Value = 1974.5454
$ws.Range("K132").Value = 5923.6362
$ws.Range("M132").Value = -3393.6362
$ws.Range("H137").Value = 2134.4119
$ws.Range("I137").Value = 2065.8333
$ws.Range("J137").Value = 2299
$ws.Range("K137").Value = 6197.499899999999
$ws.Range("L137").Value = 6897
$ws.Range("M137").Value = -3647.499899999999
$ws.Range("N137").Value = -11997
$ws.Range("N52").ClearContents()

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6453.316
$ws.Range("I61").Value = 3583.6743
$ws.Range("J61").Value = 15267.214
$ws.Range("K61").Value = 3583.6743
$ws.Range("L61").Value = 15267.214
$ws.Range("M61").Value = -3371.6743
$ws.Range("N61").Value = -15691.214
$ws.Range("H75").Value = 40173
$ws.Range("J75").Value = 40173
$ws.Range("L75").Value = 40173
$ws.Range("N75").Value = -41921
$ws.Range("H78").Value = 40173
$ws.Range("J78").Value = 40173
$ws.Range("L78").Value = 120519
$ws.Range("N78").Value = -129255
$ws.Range("H97").Value = 837.75
$ws.Range("I97").Value = 628.1429000000001
$ws.Range("J97").Value = 1326.8334
$ws.Range("K97").Value = 628.1429000000001
$ws.Range("L97").Value = 1326.8334
$ws.Range("M97").Value = -132.1429000000001
$ws.Range("N97").Value = -2318.8334
$ws.Range("H136").Value = 6453.316
$ws.Range("I136").Value = 3583.6743
$ws.Range("J136").Value = 15267.214
$ws.Range("K136").Value = 10751.0229
$ws.Range("L136").Value = 45801.642
$ws.Range("M136").Value = -8201.0229
$ws.Range("N136").Value = -50901.642

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 12000
$ws.Range("I98").Value = 12000
$ws.Range("K98").Value = 12000
$ws.Range("M98").Value = -9754
$ws.Range("H122").Value = 5807.077
$ws.Range("I122").Value = 5807.077
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 17421.231
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -14971.231
$ws.Range("N122").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3000
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("H63").Value = 3322.2
$ws.Range("I63").Value = 2406
$ws.Range("J63").Value = 3933
$ws.Range("K63").Value = 7218
$ws.Range("L63").Value = 11799
$ws.Range("M63").Value = -6469
$ws.Range("N63").Value = -13297
$ws.Range("H64").Value = 3547.1
$ws.Range("I64").Value = 2478
$ws.Range("K64").Value = 7434
$ws.Range("M64").Value = -7164
$ws.Range("H66").Value = 3322.2
$ws.Range("I66").Value = 2406
$ws.Range("J66").Value = 3933
$ws.Range("K66").Value = 21654
$ws.Range("L66").Value = 35397
$ws.Range("M66").Value = -17910
$ws.Range("N66").Value = -42885
$ws.Range("H67").Value = 3547.1
$ws.Range("I67").Value = 2478
$ws.Range("K67").Value = 7434
$ws.Range("M67").Value = -6498
$ws.Range("H68").Value = 1033.1666
$ws.Range("J68").Value = 1266.3334
$ws.Range("L68").Value = 3799.0002
$ws.Range("N68").Value = -5421.0002
$ws.Range("H71").Value = 1033.1666
$ws.Range("J71").Value = 1266.3334
$ws.Range("L71").Value = 11397.0006
$ws.Range("N71").Value = -19509.0006
$ws.Range("H107").Value = 702.6070999999999
$ws.Range("I107").Value = 297.3
$ws.Range("J107").Value = 927.7778
$ws.Range("K107").Value = 891.9000000000001
$ws.Range("L107").Value = 2783.3334
$ws.Range("M107").Value = 1028.1
$ws.Range("N107").Value = -6623.3334
$ws.Range("H113").Value = 715.625
$ws.Range("I113").Value = 716.6667
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 2150.0001
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = 19.9998999999998
$ws.Range("N113").Value = -6440
$ws.Range("H114").Value = 516.63635
$ws.Range("I114").Value = 449.6
$ws.Range("J114").Value = 572.5
$ws.Range("K114").Value = 1348.8
$ws.Range("L114").Value = 1717.5
$ws.Range("M114").Value = 1905.2
$ws.Range("N114").Value = -8225.5
$ws.Range("M58").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 10031
$ws.Range("J44").Value = 10031
$ws.Range("L44").Value = 10031
$ws.Range("N44").Value = -11223
$ws.Range("H70").Value = 5373.604
$ws.Range("I70").Value = 5298.353
$ws.Range("J70").Value = 5556.357
$ws.Range("K70").Value = 5298.353
$ws.Range("L70").Value = 5556.357
$ws.Range("M70").Value = -5028.353
$ws.Range("N70").Value = -6096.357
$ws.Range("H73").Value = 5373.604
$ws.Range("I73").Value = 5298.353
$ws.Range("J73").Value = 5556.357
$ws.Range("K73").Value = 5298.353
$ws.Range("L73").Value = 5556.357
$ws.Range("M73").Value = -4362.353
$ws.Range("N73").Value = -7428.357
$ws.Range("H75").Value = 32395.8
$ws.Range("J75").Value = 32395.8
$ws.Range("L75").Value = 32395.8
$ws.Range("N75").Value = -34143.8
$ws.Range("H78").Value = 32395.8
$ws.Range("J78").Value = 32395.8
$ws.Range("L78").Value = 97187.39999999999
$ws.Range("N78").Value = -105923.4
$ws.Range("H102").Value = 3169.6875
$ws.Range("I102").Value = 2784.44
$ws.Range("J102").Value = 4545.5713
$ws.Range("K102").Value = 2784.44
$ws.Range("L102").Value = 4545.5713
$ws.Range("M102").Value = -1162.44
$ws.Range("N102").Value = -7789.5713

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 554.8
$ws.Range("I22").Value = 350.5
$ws.Range("J22").Value = 691
$ws.Range("K22").Value = 350.5
$ws.Range("L22").Value = 691
$ws.Range("M22").Value = -55.5
$ws.Range("N22").Value = -1281
$ws.Range("H27").Value = 554.8
$ws.Range("I27").Value = 350.5
$ws.Range("J27").Value = 691
$ws.Range("K27").Value = 350.5
$ws.Range("L27").Value = 691
$ws.Range("M27").Value = -243.5
$ws.Range("N27").Value = -905
$ws.Range("H93").Value = 563.5714
$ws.Range("I93").Value = 544.5454999999999
$ws.Range("J93").Value = 633.3333
$ws.Range("K93").Value = 544.5454999999999
$ws.Range("L93").Value = 633.3333
$ws.Range("M93").Value = 703.4545000000001
$ws.Range("N93").Value = -3129.3333
$ws.Range("H100").Value = 3312.7917
$ws.Range("I100").Value = 2972.389
$ws.Range("J100").Value = 4334
$ws.Range("K100").Value = 2972.389
$ws.Range("L100").Value = 4334
$ws.Range("M100").Value = -2431.389
$ws.Range("N100").Value = -5416
$ws.Range("H122").Value = 8340.321
$ws.Range("I122").Value = 7942.5884
$ws.Range("J122").Value = 8955
$ws.Range("K122").Value = 23827.7652
$ws.Range("L122").Value = 26865
$ws.Range("M122").Value = -21377.7652
$ws.Range("N122").Value = -31765
$ws.Range("H132").Value = 3530.9614
$ws.Range("I132").Value = 2711.5557
$ws.Range("J132").Value = 5374.625
$ws.Range("K132").Value = 8134.6671
$ws.Range("L132").Value = 16123.875
$ws.Range("M132").Value = -5604.6671
$ws.Range("N132").Value = -21183.875

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 880
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 3000
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 3991.2354
$ws.Range("I132").Value = 3418.0715
$ws.Range("K132").Value = 10254.2145
$ws.Range("M132").Value = -7724.2145
$ws.Range("H136").Value = 6189.8
$ws.Range("J136").Value = 10417.941
$ws.Range("L136").Value = 31253.823
$ws.Range("N136").Value = -36353.823

Write-Output "Applied Pandaemonium_Profits.xlsx leve-profit updates"
